# "fix reframework to not use orchestrator"
# The "Settings" sheet listed Orchestrator-specific configuration
# (logF_BusinessProcessName = "Framework", OrchestratorQueueName,
# OrchestratorQueueFolder). This removes the orchestrator-queue rows and
# replaces them with P2Dispatcher-specific settings (SampleDataFolder,
# MyEmail), and renames the business-process-name value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Drop the old orchestrator-queue rows (2-5) entirely so styles/row-heights
# from the removed rows don't linger, then rebuild rows 2-4 from scratch.
$ws.Rows("2:5").Delete()

# Row 2: logF_BusinessProcessName | P2Dispatcher | Logging field...
$ws.Cells.Item(2, 1).Value = "logF_BusinessProcessName"
$ws.Cells.Item(2, 2).Value = "P2Dispatcher"
$ws.Cells.Item(2, 3).Value = "Logging field which allows grouping of log data of two or more subprocesses under the same business process name"
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Cells.Item(2, 3).WrapText = $true

# Row 3: SampleDataFolder | ..\customerqueue\ | Folder name for customer data.
$ws.Cells.Item(3, 1).Value = "SampleDataFolder"
$ws.Cells.Item(3, 2).Value = "..\customerqueue\"
$ws.Cells.Item(3, 3).Value = "Folder name for customer data."
$ws.Cells.Item(3, 3).WrapText = $true

# Row 4: MyEmail | goombaxl23xl@gmail.com | Email for sending outgoing stuff...
$ws.Cells.Item(4, 1).Value = "MyEmail"
$ws.Cells.Item(4, 2).Value = "goombaxl23xl@gmail.com"
$ws.Cells.Item(4, 3).Value = "Email for sending outgoing stuff to customers, as well as getting the initial travel plans."

# Selection moves to B6 per the saved view state.
$ws.Range("B6").Select()
